$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column (O) is being added for year 2021, to the right of the
# existing 2020 column (N). Start by copying N's formatting across into
# O for every populated row, so the new column visually matches its
# neighbour (borders, fonts, number formats, etc.).
$ws.Range("N2:N6").Copy()
$ws.Range("O2:O6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header cell: the year 2021.
$ws.Range("O3").Value = 2021

# Data rows for 2021.
$ws.Range("O6").Value = 6436.9                     # population, thousand people
$ws.Range("O5").Value = 1229.5999999999999         # waste removed, thousand tons
$ws.Range("O4").Formula = "=O5/O6*1000"            # waste per person, kg (computed)

# O5's cell uses the plain (non-rounded) numeric style already used
# elsewhere on the sheet (same look as B5/C5) rather than N5's two
# decimal rounding.
$ws.Range("B5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("O5").Value = 1229.5999999999999

# Restore the workbook's final selection.
[void]$ws.Range("P16").Select()
